# Add a new row (15) to the "Biff-Items" sheet describing a new SystemTest
# case: an item whose serial number change collides with another item that
# is also trying to claim the same serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biff-Items")

$ws.Range("A15").Value = "Z00100300022-00044"
$ws.Range("B15").Value = "CE61D09F-16BA-4A3D-90DA-2DE0D41EB06A"
$ws.Range("E15").Value = "changing to SN that another item is also trying to claim"

# Move the active selection the way the author's Excel session left it.
$ws.Range("I17").Select() | Out-Null
